$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2022" column (S) values for rows 3-33, keyed by row number.
$values = @{
  3  = 2022
  4  = 57.890663775669807
  5  = 50.022742766269019
  6  = 65.935557126462314
  7  = 41.756312988336546
  8  = 37.159844500343013
  9  = 46.265141318977122
  10 = 54.165137501740752
  11 = 47.516379220948068
  12 = 60.737168699398438
  13 = 38.800251455475774
  14 = 28.662467801253705
  15 = 48.983865238282192
  16 = 52.046737970697684
  17 = 42.039438248006412
  18 = 61.864439546842981
  19 = 54.559895023999445
  20 = 48.115700180781531
  21 = 60.917884338341217
  22 = 44.498872817808376
  23 = 40.811480640517637
  24 = 48.122142747774355
  25 = 86.286397363931727
  26 = 68.335423960459067
  27 = 104.46468013142454
  28 = 63.925654994479551
  29 = 57.785410559676791
  30 = 71.340059495655098
  31 = 50.853410128538314
  32 = 51.223021582733814
  33 = 50.500891999451071
}

# Rows that carry the bold "category" row formatting (column A uses the bold
# style on those rows) - the new S cell on these rows must be bold too,
# mirroring the rest of the row's reused font.
$boldRows = @(7, 10, 13, 16, 19, 22, 25, 28, 31)

foreach ($row in 3..33) {
    $sCell = $ws.Range("S" + $row)
    $rCell = $ws.Range("R" + $row)

    # Write the value first.
    $sCell.Value2 = $values[$row]

    # Copy R's number format / font / border so the new column matches the
    # rest of the row exactly, then (for category rows) flip to bold -
    # reusing the workbook's existing bold variant of the same font.
    $rCell.Copy()
    $sCell.PasteSpecial(-4122)  # xlPasteFormats

    if ($boldRows -contains $row) {
        $sCell.Font.Bold = $true
    }
}

$excel.CutCopyMode = 0

# Move/record the active selection like the source workbook does.
$ws.Range("T4").Select()
